# Fruta / hortaliza, semanal
# Two new weekly price records are added to the Achicoria series:
#   - one dated 2022-11-08 (serial 44873), inserted as the new row 13
#   - one dated 2022-11-10 (serial 44875), inserted as the new row 16
# All the pre-existing rows keep their relative order and simply shift
# down to make room (13,14 -> 14,15 ; 15,16,17 -> 17,18,19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the first new row at position 13 (old rows 13-17 shift to 14-18).
$ws.Rows.Item(13).Insert()

# Insert the second new row at position 16 (old rows, now at 15-18,
# from 15 down shift to 16-19).
$ws.Rows.Item(16).Insert()

# New weekly record for row 13 (2022-11-08).
$ws.Cells.Item(13,1).Value = 6
$ws.Cells.Item(13,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(13,3).Value = "Metropolitana"
$ws.Cells.Item(13,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(13,4).Value = 44873
$ws.Cells.Item(13,5).Value = 13
$ws.Cells.Item(13,6).Value = 100112010
$ws.Cells.Item(13,7).Value = "Achicoria"
$ws.Cells.Item(13,8).Value = "Sin especificar"
$ws.Cells.Item(13,9).Value = "Primera"
$ws.Cells.Item(13,10).Value = 250
$ws.Cells.Item(13,11).Value = 8000
$ws.Cells.Item(13,12).Value = 8000
$ws.Cells.Item(13,13).Value = 8000
$ws.Cells.Item(13,14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(13,15).Value = "Provincia de Quillota"
$ws.Cells.Item(13,16).Value = 500
$ws.Cells.Item(13,17).Value = 16
$ws.Cells.Item(13,18).Value = "Hortaliza"

# New weekly record for row 16 (2022-11-10).
$ws.Cells.Item(16,1).Value = 6
$ws.Cells.Item(16,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16,3).Value = "Metropolitana"
$ws.Cells.Item(16,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16,4).Value = 44875
$ws.Cells.Item(16,5).Value = 13
$ws.Cells.Item(16,6).Value = 100112010
$ws.Cells.Item(16,7).Value = "Achicoria"
$ws.Cells.Item(16,8).Value = "Sin especificar"
$ws.Cells.Item(16,9).Value = "Primera"
$ws.Cells.Item(16,10).Value = 90
$ws.Cells.Item(16,11).Value = 7000
$ws.Cells.Item(16,12).Value = 7000
$ws.Cells.Item(16,13).Value = 7000
$ws.Cells.Item(16,14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(16,15).Value = "Provincia de Quillota"
$ws.Cells.Item(16,16).Value = 438
$ws.Cells.Item(16,17).Value = 16
$ws.Cells.Item(16,18).Value = "Hortaliza"
